$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EpgScreen")

# --- Row 4: Seniour_groen_Geel -> Seniour_groen_Grijs, M4 font fixed to Regular,
#            O4 gains the "styled" look (copy format from M4, then set the value) ---
$ws.Range("A4").Value = "Seniour_groen_Grijs"
$ws.Range("M4").Value = "Proximus, ProximusRegular"
$ws.Range("M4").Copy($ws.Range("O4")) | Out-Null
$ws.Range("O4").Value = "rgba(128, 128, 128, 1)"

# --- New rows 3, 5-14: reuse the styled look of M4 (style index 3) for the M column.
#     The A-column (and B3) values are written in this particular order so that the
#     newly introduced shared strings end up laid out the same way as in the target
#     workbook. ---
$rows = @{
    6  = @{ A="Seniour_Standard_Standard";   O="rgba(255, 255, 255, 1)"; P=4 }
    9  = @{ A="Stark_Standard_Standard";     O="rgba(255, 255, 255, 1)"; P=6 }
    3  = @{ A="Seniour_groen_Standard";      O="rgba(255, 255, 255, 1)"; P=4 }
    10 = @{ A="Strak_groen_Standard";        O="rgba(255, 255, 255, 1)"; P=6 }
    11 = @{ A="Strak_groen_grijs";           O="rgba(128, 128, 128, 1)"; P=6 }
    12 = @{ A="Strak_groen_geel";            O="rgba(255, 255, 0, 1)";   P=6 }
    13 = @{ A="Strak_Standard_grijs";        O="rgba(128, 128, 128, 1)"; P=6 }
    14 = @{ A="Strak_Standard_geel";         O="rgba(255, 255, 0, 1)";   P=6 }
    7  = @{ A="Seniour_Standard_geel";       O="rgba(255, 255, 0, 1)";   P=4 }
    8  = @{ A="Seniour_Standard_grijs";      O="rgba(128, 128, 128, 1)"; P=4 }
    5  = @{ A="Seniour_groen_Geel";          O="rgba(255, 255, 0, 1)";   P=4 }
}

$order = @(6, 9, 3, 10, 11, 12, 13, 14, 7, 8, 5)
foreach ($n in $order) {
    $r = $rows[$n]
    $ws.Range("A$n").Value = $r.A
    $ws.Range("B$n").Value = "programTitle"
    $ws.Range("L$n").Value = "24px"
    $ws.Range("M4").Copy($ws.Range("M$n")) | Out-Null
    $ws.Range("M$n").Value = "Proximus, ProximusRegular"
    $ws.Range("O$n").Value = $r.O
    $ws.Range("P$n").Value = $r.P
}

# B3 keeps the trailing-space "programTitle " text, and O3 reverts to the
# un-styled look (ClearFormats removes the inherited style index).
$ws.Range("B3").Value = "programTitle "
$ws.Range("O3").Value = "rgba(255, 255, 255, 1)"
$ws.Range("O3").ClearFormats() | Out-Null

# --- Row 17 / 18: a couple of stray single-space cells ---
$ws.Range("A17").Value = " "
$ws.Range("M8").Copy($ws.Range("M18")) | Out-Null
$ws.Range("M18").Value = " "

# --- Update the active selection shown when the workbook is opened ---
$ws.Range("B6").Select() | Out-Null
